$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# realistic_candidate_capacities_tobe_installed -> FALSE
$ws.Range("B19").Value = $false

# realistic_candidate_capacities_to_test -> FALSE (drives C20 text formula result)
$ws.Range("B20").Value = $false

# start_dismantling_tick -> formula referencing B9 instead of hard-coded 100
$ws.Range("B23").Formula = "=B9"
